# "reattempt snap shot array and first bad version"
#
# - First Bad Version (row 8): bump Python attempt rating from 4* to 5*
# - Snapshot Array (row 15): bump Python attempt rating from 2 to 3
# - Maximum Depth of Binary Tree (Recursive) (row 24): bump Python attempts
#   1 -> 2, and fill in Time Complexity as O(n)
# - Maximum Depth of Binary Tree (BFS) (row 26): fill in Time Complexity as O(n)
# - Break a palindrome (row 28): re-enter Python attempts value so it picks
#   up the normal data-row formatting
# - Minor row height touch-ups picked up while reviewing the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First Bad Version: Python attempts 4* -> 5*
$ws.Range("E8").Value = "5*"

# Snapshot Array: Python attempts 2 -> 3
$ws.Range("E15").Value = 3

# Maximum Depth of Binary Tree (Recursive): Python attempts 1 -> 2, T Complexity -> O(n)
$ws.Range("E24").Value = 2
$ws.Range("K24").Value = "O(n)"

# Maximum Depth of Binary Tree (BFS): T Complexity -> O(n)
$ws.Range("K26").Value = "O(n)"

# Break a palindrome: re-key the Python attempts cell so it matches the
# standard data-row style used elsewhere in the column (same format as the
# other populated "Python attempts" cells, e.g. E24)
$ws.Range("E28").Value = 1
$ws.Range("E24").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row-height touch ups
$ws.Rows(16).RowHeight = 44.25
$ws.Rows(31).RowHeight = 19.5
$ws.Rows(32).RowHeight = 19.5
$ws.Rows(33).RowHeight = 19.5
$ws.Rows(34).RowHeight = 19.5
